$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6 and 7 (now-removed data rows) first so remaining rows shift up correctly
$ws.Rows.Item(7).Delete() | Out-Null
$ws.Rows.Item(6).Delete() | Out-Null

# Row 2
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Gdnf"
$ws.Range("C2").Value = "Gfra2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.509729
$ws.Range("H2").Value = 1.529187
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5298676666666666
$ws.Range("N2").Value = 1.589603
$ws.Range("O2").Value = 0.0918795741734264
$ws.Range("P2").Value = 0.09187957417342642
$ws.Range("Q2").Value = 0.2700889158623333
$ws.Range("R2").Value = 2.430800242760999
$ws.Range("S2").Value = 0.0918795741734264
$ws.Range("T2").Value = 0.09187957417342642

# Row 3
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Gdnf"
$ws.Range("C3").Value = "Gfra2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.509729
$ws.Range("H3").Value = 1.529187
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.390723333333333
$ws.Range("N3").Value = 7.17217
$ws.Range("O3").Value = 0.4145537756908006
$ws.Range("P3").Value = 0.4145537756908007
$ws.Range("Q3").Value = 1.218621013976667
$ws.Range("R3").Value = 10.96758912579
$ws.Range("S3").Value = 0.4145537756908006
$ws.Range("T3").Value = 0.4145537756908007

# Row 4
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Gdnf"
$ws.Range("C4").Value = "Gfra2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.509729
$ws.Range("H4").Value = 1.529187
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.463531333333334
$ws.Range("N4").Value = 7.390594
$ws.Range("O4").Value = 0.4271787544491802
$ws.Range("P4").Value = 0.4271787544491803
$ws.Range("Q4").Value = 1.255733363008667
$ws.Range("R4").Value = 11.301600267078
$ws.Range("S4").Value = 0.4271787544491802
$ws.Range("T4").Value = 0.4271787544491803

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Gdnf"
$ws.Range("C5").Value = "Gfra2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.509729
$ws.Range("H5").Value = 1.529187
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3828576666666666
$ws.Range("N5").Value = 1.148573
$ws.Range("O5").Value = 0.06638789568659274
$ws.Range("P5").Value = 0.06638789568659276
$ws.Range("Q5").Value = 0.1951536555723333
$ws.Range("R5").Value = 1.756382900151
$ws.Range("S5").Value = 0.06638789568659274
$ws.Range("T5").Value = 0.06638789568659276
